$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting of column I (years 2020 column) into the two new
# columns J (2021) and K (2022) for the data rows (header row 4 through the
# total row 14), then fill in the new values.
$ws.Range("I4:I14").Copy() | Out-Null
$ws.Range("J4").PasteSpecial(-4122) | Out-Null
$ws.Range("I4:I14").Copy() | Out-Null
$ws.Range("K4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# New year headers
$ws.Range("J4").Value = 2021
$ws.Range("K4").Value = 2022

# New data values for rows 5-14
$ws.Range("J5").Value = 34.075233127500141
$ws.Range("K5").Value = 35.305353068702679

$ws.Range("J6").Value = 44.487602536118636
$ws.Range("K6").Value = 49.31549563692068

$ws.Range("J7").Value = 40.668697007891453
$ws.Range("K7").Value = 45.444207273635158

$ws.Range("J8").Value = 50.797011639929529
$ws.Range("K8").Value = 46.810603774236895

$ws.Range("J9").Value = 46.848562449074493
$ws.Range("K9").Value = 45.450816127137941

$ws.Range("J10").Value = 44.458036086558309
$ws.Range("K10").Value = 50.474514452886076

$ws.Range("J11").Value = 40.532201616746903
$ws.Range("K11").Value = 40.14796186663478

$ws.Range("J12").Value = 33.353175884696697
$ws.Range("K12").Value = 29.735683954543184

$ws.Range("J13").Value = 10.46405303463253
$ws.Range("K13").Value = 12.912087912087852

$ws.Range("J14").Value = 43.479082661290349
$ws.Range("K14").Value = 41.117034465658314

# Update the selection to match the post-edit state
$ws.Range("M6").Select()
